$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.164.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.74%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.209.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.35%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.65%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.56%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.55%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.19%  '

# Row 10
$ws.Range("E10").Value = '  +4.41%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.437'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.79%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.768.77'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.65%  '

# Row 13
$ws.Range("E13").Value = '  -1.23%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.81%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000175'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.73%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.265.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.74%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.243.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.55%  '

# Row 18
$ws.Range("E18").Value = '  +0.66%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.33%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.86%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.86%  '

# Row 22
$ws.Range("E22").Value = '  +0.52%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.530'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.89%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.96%  '

# Row 26
$ws.Range("E26").Value = '  +2.19%  '

# Row 27
$ws.Range("E27").Value = '  +0.08%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0913'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.02%  '

# Row 29
$ws.Range("E29").Value = '  +0.85%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.21%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.75%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.61%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.88%  '

# Row 34
$ws.Range("E34").Value = '  +6.44%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '157.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.70%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.37'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.01%  '

# Row 37
$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.99'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.28%  '

# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.793.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.75%  '

# Row 39
$ws.Range("E39").Value = '  +4.90%  '

# Row 40
$ws.Range("E40").Value = '  +0.75%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.33%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.721'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.06%  '

# Row 44
$ws.Range("E44").Value = '  +4.33%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.255.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.50%  '

# Row 46
$ws.Range("E46").Value = '  +3.47%  '

# Row 47
$ws.Range("E47").Value = '  +0.08%  '

# Row 48
$ws.Range("E48").Value = '  -0.45%  '

# Row 49
$ws.Range("E49").Value = '  +7.79%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.11%  '

# Row 51
$ws.Range("E51").Value = '  +0.02%  '
